# Apply the "encode" edit: introduce typo'd header names and spell out the
# y/n flag column as yes/no, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1): rename columns, introducing the typo'd labels ---
$ws.Range("A1").Value = "accessibilityes"
$ws.Range("B1").Value = "defenosive"
$ws.Range("C1").Value = "noum_forces"
$ws.Range("D1").Value = "visibilityes"
$ws.Range("E1").Value = "water_level"
$ws.Range("F1").Value = "soil_tyespe"
$ws.Range("G1").Value = "topographyes"
$ws.Range("H1").Value = "speed"
$ws.Range("I1").Value = "importanoce"
$ws.Range("J1").Value = "output"

# Header row grows taller (wraps across two lines in Excel after rename)
$ws.Rows.Item(1).RowHeight = 32.25

# --- Column D ("visibility") data rows: spell "y"/"n" out as "yes"/"no" ---
for ($r = 2; $r -le 24; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $v = $cell.Value()
    if ($v -eq "y") {
        $cell.Value = "yes"
    } elseif ($v -eq "n") {
        $cell.Value = "no"
    }
}

# --- View state: scroll position / active selection moved ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H17").Select() | Out-Null
